$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @("65.803.69", "  -0.51%  ")
    3 = @("3.318.14", "  +1.21%  ")
    4 = @("0.999", "  -0.08%  ")
    5 = @("557.29", "  -0.21%  ")
    6 = @("184.94", "  -0.94%  ")
    7 = @($null, "  +0.17%  ")
    8 = @("3.310.32", "  +1.18%  ")
    9 = @("0.575", "  -2.90%  ")
    10 = @($null, "  -7.46%  ")
    11 = @("0.576", "  -1.92%  ")
    12 = @($null, "  -3.88%  ")
    13 = @("0.0000261", "  -2.26%  ")
    14 = @("3.849.86", "  +1.27%  ")
    15 = @("8.42", "  -2.69%  ")
    16 = @("568.94", "  -10.44%  ")
    17 = @("65.835.61", "  -0.25%  ")
    19 = @("3.317.22", "  +1.26%  ")
    20 = @($null, "  -1.64%  ")
    21 = @("10.80", "  -5.09%  ")
    22 = @($null, "  -2.02%  ")
    23 = @("17.95", "  -2.71%  ")
    24 = @("4.98", "  +0.94%  ")
    25 = @("97.62", "  -8.92%  ")
    26 = @("3.93", "  -0.97%  ")
    27 = @("2.68", "  +0.15%  ")
    28 = @("9.34", "  -2.83%  ")
    29 = @("8.44", "  -3.38%  ")
    30 = @($null, "  +0.07%  ")
    31 = @("6.66", "  +6.01%  ")
    32 = @($null, "  -9.86%  ")
    33 = @("559.93", "  +4.97%  ")
    34 = @("10.80", "  -2.33%  ")
    35 = @($null, "  -2.29%  ")
    36 = @("3.733.77", "  +0.05%  ")
    37 = @("0.999", "  -0.15%  ")
    38 = @("55.51", "  -3.52%  ")
    39 = @("33.61", "  +2.14%  ")
    40 = @($null, "  -4.42%  ")
    41 = @("0.0₃0682", "  -7.14%  ")
    42 = @($null, "  -8.13%  ")
    43 = @($null, "  -6.40%  ")
    44 = @($null, "  +1.47%  ")
    45 = @($null, "  -2.23%  ")
    46 = @($null, "  -2.07%  ")
    47 = @("2.98", "  -13.60%  ")
    48 = @($null, "  -1.96%  ")
    49 = @($null, "  +0.20%  ")
    50 = @($null, "  -4.13%  ")
    51 = @("124.69", "  +1.91%  ")
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($dVal -ne $null) {
        $dCell = $ws.Cells.Item($r, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
        $dCell.Style = "Normal"
    }
    $ws.Cells.Item($r, 5).Value = $eVal
}

Write-Host "Updated $($data.Keys.Count) rows"
